# The author regenerated the controlled-vocabulary sheet from an updated
# Google Sheet export. The new export drops the old row 512 ("cl:10583" /
# "catastrophic event") entirely, which causes every subsequent row
# (old rows 513-680) to shift up by one (new rows 512-679). The last
# row of the old sheet (680, "phenology") therefore disappears and the
# sheet's used range shrinks from A1:BE680 to A1:BE679.
#
# Deleting the entire row 512 reproduces exactly that: Excel removes the
# row, shifts everything below it up by one, and recalculates the sheet
# dimension automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(512).Delete()
